$d = $word.ActiveDocument

# Locate the paragraph that contains the "{{response_4}}" placeholder text.
$paras = $d.Paragraphs
$n = $paras.Count
$targetIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*{{response_4}}*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    # Remove the empty (italic / rtl) paragraph that immediately precedes the
    # placeholder paragraph, including its paragraph mark.
    $prev = $paras.Item($targetIndex - 1)
    $prev.Range.Delete()
}

# Clear the "{{response_4}}" placeholder text, leaving the (now empty)
# paragraph mark / pPr in place.
$d.Content.Find.Execute("{{response_4}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
